$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly data refresh: a new observation (Primera/Segunda) is inserted at the
# top of the historical block (row 107), pushing all existing rows down by
# two and extending the sheet from R220 to R222.
$ws.Rows.Item(107).Insert()
$ws.Rows.Item(107).Insert()

# New row 107 - Primera
$ws.Cells.Item(107,1).Value = 1
$ws.Cells.Item(107,2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(107,3).Value = "Arica y Parinacota"
$ws.Cells.Item(107,4).Value = 44494
$ws.Cells.Item(107,5).Value = 15
$ws.Cells.Item(107,6).Value = 100112032
$ws.Cells.Item(107,7).Value = "Zapallo italiano"
$ws.Cells.Item(107,8).Value = "Huracán"
$ws.Cells.Item(107,9).Value = "Primera"
$ws.Cells.Item(107,10).Value = 120
$ws.Cells.Item(107,11).Value = 7000
$ws.Cells.Item(107,12).Value = 8000
$ws.Cells.Item(107,13).Value = 7500
$ws.Cells.Item(107,14).Value = "`$/caja 70 unidades"
$ws.Cells.Item(107,15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(107,16).Value = 107
$ws.Cells.Item(107,17).Value = 70
$ws.Cells.Item(107,18).Value = "Hortaliza"

# New row 108 - Segunda
$ws.Cells.Item(108,1).Value = 1
$ws.Cells.Item(108,2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(108,3).Value = "Arica y Parinacota"
$ws.Cells.Item(108,4).Value = 44494
$ws.Cells.Item(108,5).Value = 15
$ws.Cells.Item(108,6).Value = 100112032
$ws.Cells.Item(108,7).Value = "Zapallo italiano"
$ws.Cells.Item(108,8).Value = "Huracán"
$ws.Cells.Item(108,9).Value = "Segunda"
$ws.Cells.Item(108,10).Value = 130
$ws.Cells.Item(108,11).Value = 5000
$ws.Cells.Item(108,12).Value = 6000
$ws.Cells.Item(108,13).Value = 5500
$ws.Cells.Item(108,14).Value = "`$/caja 100 unidades"
$ws.Cells.Item(108,15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(108,16).Value = 55
$ws.Cells.Item(108,17).Value = 100
$ws.Cells.Item(108,18).Value = "Hortaliza"
